# Update imputed values in result_data_RandomForest.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = -12.9215
$ws.Range("A12").Value = -21.36689999999999
$ws.Range("C23").Value = -12.2058
$ws.Range("C28").Value = -13.1473
$ws.Range("A32").Value = -21.5947
$ws.Range("C32").Value = -12.88529999999999
$ws.Range("C34").Value = -12.08330000000001
$ws.Range("A36").Value = -20.1876
$ws.Range("A38").Value = -19.70719999999999
$ws.Range("C42").Value = -12.74870000000001
$ws.Range("A46").Value = -21.93250000000001
$ws.Range("A54").Value = -22.1041
$ws.Range("C54").Value = -12.2146
$ws.Range("A55").Value = -22.2074
$ws.Range("A67").Value = -21.41939999999997
$ws.Range("A69").Value = -21.66559999999997
$ws.Range("A72").Value = -21.8633
$ws.Range("A91").Value = -20.54029999999998
$ws.Range("C97").Value = -11.94760000000001
$ws.Range("A99").Value = -22.0018
$ws.Range("C99").Value = -13.1149
$ws.Range("C101").Value = -13.055
$ws.Range("A104").Value = -21.34149999999999
